$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Remove the stray "_GoBack" bookmark that currently sits between
#    "...and so on. " and "Each time we encountered..." and merge
#    those two runs back into a single run (matches a normal Word
#    edit that touched that span and coalesced the run).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

$mergeText1 = "so on. Each time we encountered"
$d.Content.Find.Execute($mergeText1, $true, $false, $false, $false, $false,
                         $true, 1, $false, $mergeText1, 2)

# ------------------------------------------------------------------
# 2) Fix "github" -> "GitHub" and drop the spellcheck proofErr
#    wrapper that surrounded the lowercase word. We do this by first
#    merging the run across both of its neighbours (which also
#    swallows the two <w:proofErr/> markers), then re-isolating the
#    "GitHub" word back into its own run using a pair of bookmarks
#    (bookmarks always sit between runs, forcing a split) - the
#    temporary one is removed afterwards, leaving only the final
#    "_GoBack" bookmark right after "GitHub", matching a real edit
#    session's last-edit marker.
# ------------------------------------------------------------------
$d.Content.Find.Execute("use github", $true, $false, $false, $false, $false,
                         $true, 1, $false, "use GitHub", 2)

$d.Content.Find.Execute("GitHub to", $true, $false, $false, $false, $false,
                         $true, 1, $false, "GitHub to", 2)

$gh = $d.Content
$gh.Find.Execute("GitHub")
$ghStart = $gh.Start
$ghEnd = $gh.End

$beforeGh = $d.Range($ghStart, $ghStart)
$d.Bookmarks.Add("_TmpSplit", $beforeGh)

$afterGh = $d.Range($ghEnd, $ghEnd)
$d.Bookmarks.Add("_GoBack", $afterGh)

$d.Bookmarks.Item("_TmpSplit").Delete()

# ------------------------------------------------------------------
# 3) Merge the lone-space run with the following "The project was
#    scoped..." run into a single run (keeps the separate run that
#    ends in "...have been created." intact).
# ------------------------------------------------------------------
$mergeText2 = "project was scoped well"
$d.Content.Find.Execute($mergeText2, $true, $false, $false, $false, $false,
                         $true, 1, $false, $mergeText2, 2)
